# docs: fix sample questionnaire (Fixes #9248)
#
# Adds a "packages" worksheet (name/label/description) in front of the
# existing "entities" and "attributes" sheets, wires the new "doc" package
# into "entities" (new "package" column) and into "attributes" (entity ids
# renamed to "doc_SimpleQuestionnaire", new "labelAttribute" column, a new
# "name" attribute row, and a corrected drivers-license question label).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "packages" sheet in front of everything else.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "packages"
$newSheet.Move($wb.Worksheets.Item(1))

# NB: after .Move() the old $newSheet COM reference no longer tracks the
# moved sheet (it rebinds by position) - always re-fetch by index/name.
$packages = $wb.Worksheets.Item(1)

$packages.Cells.Item(1, 1).Value = "name"
$packages.Cells.Item(1, 2).Value = "label"
$packages.Cells.Item(1, 3).Value = "description"

$packages.Cells.Item(2, 1).Value = "doc"
$packages.Cells.Item(2, 2).Value = "Documentation"
$packages.Cells.Item(2, 3).Value = "Sample data"

[void]$packages.Range("D1:E1048576").Select()

# ---------------------------------------------------------------------
# 2. "entities" sheet: insert a "package" column right after "name".
# ---------------------------------------------------------------------
$entities = $wb.Worksheets.Item(2)

$entities.Columns.Item(2).Insert()
$entities.Cells.Item(1, 2).Value = "package"
$entities.Cells.Item(2, 2).Value = "doc"

[void]$entities.Range("B2").Select()

# ---------------------------------------------------------------------
# 3. "attributes" sheet: add "labelAttribute" column, add a "name"
#    attribute row, point entities at the new "doc_SimpleQuestionnaire"
#    id, and fix up the drivers-license question text.
# ---------------------------------------------------------------------
$attributes = $wb.Worksheets.Item(3)

# New column F = labelAttribute (pushes nillable/visible/label/description
# one column to the right, carrying their styles - incl. the highlighted
# formula cell - along with them).
$attributes.Columns.Item(6).Insert()
$attributes.Cells.Item(1, 6).Value = "labelAttribute"

# New row 3 = the "name" attribute (pushes age/driverslicence down one row).
$attributes.Rows.Item(3).Insert()

$attributes.Cells.Item(3, 1).Value = "doc_SimpleQuestionnaire"
$attributes.Cells.Item(3, 2).Value = "name"
$attributes.Cells.Item(3, 3).Value = "string"
$attributes.Cells.Item(3, 6).Value = $true
$attributes.Cells.Item(3, 7).Value = $false
$attributes.Cells.Item(3, 9).Value = "What is your name?"
$attributes.Cells.Item(3, 10).Value = "Name"

# Existing rows now point at the "doc_SimpleQuestionnaire" entity id.
$attributes.Cells.Item(2, 1).Value = "doc_SimpleQuestionnaire"
$attributes.Cells.Item(4, 1).Value = "doc_SimpleQuestionnaire"
$attributes.Cells.Item(5, 1).Value = "doc_SimpleQuestionnaire"

# Drivers-license question label gains a trailing question mark.
$attributes.Cells.Item(5, 9).Value = "Do you have a drivers license?"

[void]$attributes.Range("I6").Select()
$attributes.Activate()
